$p = $ppt.ActivePresentation

# Slide 14: "This is a complete 180 to how DAM is implemented at Jostens
# and is a return to the past" -> reworded first paragraph, split across
# several runs (as PowerPoint does when text is edited/retyped in pieces).
$s = $p.Slides.Item(14)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

$oldFirstPara = "This is a complete 180 to how DAM is implemented at Jostens and is a return to the past"
$newFirstPara = "This is a pretty radical shift to how DAM is implemented at Jostens with MediaBin and is a return to the past"

# Replace just the paragraph's characters (not the trailing paragraph mark).
$firstPara = $tr.Characters(1, $oldFirstPara.Length)
$firstPara.Text = $newFirstPara

# Re-set each new segment's own text (identical content) so the single run
# produced above gets split into separate runs at the right boundaries,
# matching the six <a:r> runs in the target markup.
$seg1Len = "This is a ".Length
$seg2Len = "pretty radical shift to ".Length
$seg3Len = "how DAM is implemented at Jostens ".Length
$seg4Len = "with MediaBin ".Length
$seg5Len = "and ".Length
$seg6Len = "is a return to the past".Length

$seg2Start = $seg1Len + 1
$seg2 = $tr.Characters($seg2Start, $seg2Len)
$seg2.Text = $seg2.Text

$seg4Start = $seg2Start + $seg2Len + $seg3Len
$seg4 = $tr.Characters($seg4Start, $seg4Len)
$seg4.Text = $seg4.Text

$seg5Start = $seg4Start + $seg4Len
$seg5 = $tr.Characters($seg5Start, $seg5Len)
$seg5.Text = $seg5.Text

$seg6Start = $seg5Start + $seg5Len
$seg6 = $tr.Characters($seg6Start, $seg6Len)
$seg6.Text = $seg6.Text
